# The deck's theme parts (ppt/theme/theme1.xml "Office Theme" and
# ppt/theme/theme2.xml "Integral") had their contents swapped: the slide
# master's theme (the one exposed through the PowerPoint object model as
# Design / Master.Theme / ThemeColorScheme) switches from the "Integral"
# palette to the original default "Office" palette.
#
# Helper: turn an "RRGGBB" hex string into the BGR-packed integer that the
# ColorFormat.RGB COM property expects (PowerPoint stores RGB as 0x00BBGGRR).
function ConvertTo-BgrInt($hex) {
  $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
  $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
  $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
  return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$scheme = $master.Theme.ThemeColorScheme

# Target palette = the original Office theme colors (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - standard msoThemeColorIndex order 1..12).
$officeColors = @(
  "000000",
  "FFFFFF",
  "44546A",
  "E7E6E6",
  "5B9BD5",
  "ED7D31",
  "A5A5A5",
  "FFC000",
  "4472C4",
  "70AD47",
  "0563C1",
  "954F72"
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
  $scheme.Item($i).RGB = ConvertTo-BgrInt $officeColors[$i - 1]
}
